# edit.ps1
# Applies two kinds of changes described by the diff:
#   1. Renames several worksheet tabs to add spaces / ampersands for readability.
#   2. Bumps the "days" component of every "Age" column (column E) value
#      (formatted as "YY-DDD", years-days) forward by one day on every
#      per-player stats sheet (i.e. every sheet except "Matches").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename sheets
# ---------------------------------------------------------------------------
$renames = @{
    "StandardStats"    = "Standard Stats"
    "ShootingStats"    = "Shooting Stats"
    "PassingStats"     = "Passing Stats"
    "PassTypes"        = "Pass Types"
    "GoalShotCreation" = "Goal & Shot Creation"
    "DefensiveActions" = "Defensive Actions"
    "PlayingTime"      = "Playing Time"
    "MiscStats"        = "Miscellaneous Stats"
}

foreach ($oldName in $renames.Keys) {
    $sheet = $wb.Worksheets.Item($oldName)
    $sheet.Name = $renames[$oldName]
}

# ---------------------------------------------------------------------------
# 2. Bump the "Age" column (column E) day-of-year counter by one for every
#    data row on every stats sheet (everything except "Matches").
# ---------------------------------------------------------------------------
function Bump-AgeDays {
    param($val)
    $parts = $val -split '-'
    $years = $parts[0]
    $days = [int]$parts[1]
    $newDays = $days + 1
    return ("{0}-{1:D3}" -f $years, $newDays)
}

$statSheetNames = @(
    "Standard Stats",
    "Shooting Stats",
    "Passing Stats",
    "Pass Types",
    "Goal & Shot Creation",
    "Defensive Actions",
    "Possession",
    "Playing Time",
    "Miscellaneous Stats"
)

foreach ($sheetName in $statSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    $lastRow = $used.Row + $used.Rows.Count - 1

    for ($r = 4; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 5)
        $oldVal = $cell.Value2()
        if ($oldVal -and $oldVal -match '^\d+-\d+$') {
            $cell.Value2 = Bump-AgeDays $oldVal
        }
    }
}
